$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a "last changed" date serial value that was
# bumped by one day (2023-09-02 -> 2023-09-03, serial 45171 -> 45172) for
# every data row (rows 2 through 506).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 506) { $lastRow = 506 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45172
